$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos data refresh (prices + volume deltas), some rows also swap rank-neighbour identity

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.516.03'
$ws.Range("E2").Value = '  -1.16%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.690.74'
$ws.Range("E3").Value = '  -0.76%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.42%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.20'
$ws.Range("E5").Value = '  -0.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.32%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3911'
$ws.Range("E7").Value = '  -0.83%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4042'
$ws.Range("E8").Value = '  -0.16%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.491'
$ws.Range("E9").Value = '  -1.81%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9991'
$ws.Range("E10").Value = '  -0.46%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.17'
$ws.Range("E11").Value = '  -0.54%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08790'
$ws.Range("E12").Value = '  -1.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.47'
$ws.Range("E13").Value = '  +11.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.482'
$ws.Range("E14").Value = '  +1.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.144'
$ws.Range("E15").Value = '  +1.80%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001349'
$ws.Range("E16").Value = '  +1.61%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.687.58'
$ws.Range("E17").Value = '  -1.06%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '98.20'
$ws.Range("E18").Value = '  -1.97%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07243'
$ws.Range("E19").Value = '  +2.83%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.22'
$ws.Range("E20").Value = '  +2.33%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.320'
$ws.Range("E21").Value = '  +3.19%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.31'
$ws.Range("E23").Value = '  -1.04%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.499.48'
$ws.Range("E24").Value = '  -1.15%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.028'
$ws.Range("E25").Value = '  -5.96%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.339'
$ws.Range("E26").Value = '  -1.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.65'
$ws.Range("E27").Value = '  -0.64%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '167.55'
$ws.Range("E28").Value = '  +3.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.443'
$ws.Range("E29").Value = '  -2.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.414'
$ws.Range("E30").Value = '  +4.56%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '138.52'
$ws.Range("E31").Value = '  +1.33%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.872.99'
$ws.Range("E32").Value = '  -1.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.08725'
$ws.Range("E33").Value = '  -1.37%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.316'
$ws.Range("E34").Value = '  -3.75%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.113'
$ws.Range("E35").Value = '  +6.29%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.042'
$ws.Range("E36").Value = '  -4.19%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02985'
$ws.Range("E37").Value = '  +7.32%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2787'
$ws.Range("E38").Value = '  +0.87%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.87'
$ws.Range("E39").Value = '  -2.15%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8080'
$ws.Range("E40").Value = '  +4.54%  '

# Row 41
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09149'
$ws.Range("E41").Value = '  -0.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.20'
$ws.Range("E42").Value = '  -3.17%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.476'
$ws.Range("E43").Value = '  +1.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.51'
$ws.Range("E44").Value = '  +9.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.671'
$ws.Range("E45").Value = '  +3.85%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7267'
$ws.Range("E46").Value = '  +0.59%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.260'
$ws.Range("E47").Value = '  +1.20%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.407'
$ws.Range("E48").Value = '  +5.66%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9993'
$ws.Range("E49").Value = '  -0.29%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '139.69'
$ws.Range("E50").Value = '  -0.89%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08126'
$ws.Range("E51").Value = '  +1.44%  '

